# Add a new Check Link row: "Consensi Email sui Contratti" under the
# Clients (Anagrafe) section of the BurgerMenu test plan sheet.
#
# This mirrors what happens when, in Excel, the author right-clicks row 13
# and chooses "Insert" (shifting the existing rows 13..90 down to 14..91)
# and then fills in the new row with data matching the surrounding
# "Clients" rows (rows 2-12), and finally re-points the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$newRow = 13

# Insert a blank row, pushing current rows 13..90 down to 14..91.
$ws.Rows.Item($newRow).Insert()

# Bring over the formatting (styles / wrap text / borders) used by the
# surrounding "Clients" rows so the new row looks consistent with its
# neighbours (same as rows 2-12 use: s=3 for col A, s=4 for cols B-J).
$ws.Range("A14:J14").Copy()
$ws.Range("A13:J13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Match the row height used by the rest of this block (75pt).
$ws.Rows.Item($newRow).RowHeight = 75

# Fill in the new Check Link data.
$ws.Range("A13").Value = "MatrixWeb: Navigation BurgerMenuClients_Verifica aggancio Consensi Email sui Contratti"
$ws.Range("B13").Value = "Verifica aggancio Consensi Email sui Contratti"
$ws.Range("C13").Value = "Si accede a Clients, click burgerMenu e verifica atterraggio della pagina"
$ws.Range("D13").Value = "Pusateri Kevin (Leased Employed)"
$ws.Range("E13").Value = "Design"
$ws.Range("F13").Value = "Planned"
$ws.Range("G13").Value = "Anagrafe"
$ws.Range("H13").Value = "Automation"
$ws.Range("I13").Value = "Test Factory"
$ws.Range("J13").Value = "Allianz Projects\Digital Interaction\Allianz Matrix Web\Clients"

# Reflect the cursor / selection position recorded for this edit.
$ws.Range("B7").Select()
$ws.Range("J13").Select()
